$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto market data (prices / 1h volume %) matching the
# GitHub Actions data refresh commit.

$ws.Range('D2').Value = '62.307.13'
$ws.Range('E2').Value = '  +1.24%  '
$ws.Range('D3').Value = '3.428.65'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '406.76'
$ws.Range('E5').Value = '  +0.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '130.68'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.595'
$ws.Range('E7').Value = '  -2.63%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.691'
$ws.Range('E9').Value = '  +2.09%  '
$ws.Range('E10').Value = '  +7.40%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '41.91'
$ws.Range('E11').Value = '  -1.22%  '
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '19.83'
$ws.Range('E13').Value = '  +0.42%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '8.36'
$ws.Range('E14').Value = '  -2.17%  '
$ws.Range('D15').Value = '3.412.27'
$ws.Range('E15').Value = '  +0.58%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '62.255.74'
$ws.Range('E16').Value = '  +1.43%  '
$ws.Range('B17').Value = 'Uniswap'
$ws.Range('C17').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '11.56'
$ws.Range('E17').Value = '  +0.77%  '
$ws.Range('E18').Value = '  -1.11%  '
$ws.Range('E19').Value = '  +11.73%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.16'
$ws.Range('E20').Value = '  -3.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '84.02'
$ws.Range('E21').Value = '  +1.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '312.15'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '12.79'
$ws.Range('E23').Value = '  -2.77%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '3.16'
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '4.73'
$ws.Range('E25').Value = '  +1.12%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.63'
$ws.Range('E26').Value = '  -0.46%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.13'
$ws.Range('E27').Value = '  -5.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.84'
$ws.Range('E28').Value = '  +4.63%  '
$ws.Range('E29').Value = '  +5.47%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '44.74'
$ws.Range('E30').Value = '  +4.99%  '
$ws.Range('E31').Value = '  -1.15%  '
$ws.Range('E32').Value = '  -0.89%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '11.34'
$ws.Range('E33').Value = '  -4.13%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('E35').Value = '  -0.48%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '51.83'
$ws.Range('E36').Value = '  -0.95%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.999'
$ws.Range('E37').Value = '  +0.24%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.98'
$ws.Range('E38').Value = '  -1.26%  '
$ws.Range('E39').Value = '  +12.28%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.30'
$ws.Range('E40').Value = '  -4.10%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '142.16'
$ws.Range('E41').Value = '  +3.72%  '
$ws.Range('E42').Value = '  -0.35%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.97'
$ws.Range('E43').Value = '  -3.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.90'
$ws.Range('E44').Value = '  -3.42%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '16.78'
$ws.Range('E45').Value = '  -1.67%  '
$ws.Range('E46').Value = '  -0.05%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '21.27'
$ws.Range('E47').Value = '  -2.40%  '
$ws.Range('D48').Value = '2.101.53'
$ws.Range('E48').Value = '  -2.34%  '
$ws.Range('B49').Value = 'ThetaToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.97'
$ws.Range('E49').Value = '  +2.54%  '
$ws.Range('B50').Value = 'ApeXProtocol'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.31'
$ws.Range('E50').Value = '  -1.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.09'
$ws.Range('E51').Value = '  +26.98%  '
